$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last four rows of the "Approved/Rejected" column (I9:I12) are being
# changed from "Approved" to "Rejected".
$ws.Range("I9:I12").Value = "Rejected"

# Reflect the newly-edited range as the current selection (active cell I9,
# selected range I9:I12) just like the author left it.
$excel.Goto($ws.Range("I9:I12"))
